$d = $word.ActiveDocument

# 1) "ihmisille mahdollisuus hyväksyä kuljetuksia, vastaanottaa semmosia"
#    was split across runs with proofErr spell-check markers around
#    "semmosia". Rebuild the paragraph as a single clean run.
$p1 = $d.Paragraphs(3)
$r1 = $p1.Range
$start1 = $r1.Start
$r1.Delete()
$ins1 = $d.Range($start1, $start1)
$ins1.InsertBefore("ihmisille mahdollisuus hyväksyä kuljetuksia, vastaanottaa semmosia`r")

# 2) "tori? yms" was likewise split across runs with proofErr markers
#    around "yms". Rebuild the paragraph as a single clean run.
$p2 = $d.Paragraphs(4)
$r2 = $p2.Range
$start2 = $r2.Start
$r2.Delete()
$ins2 = $d.Range($start2, $start2)
$ins2.InsertBefore("tori? yms`r")

# 3) Insert a new paragraph "sikapestään" right after "kilometrikorvaus?"
$p3 = $d.Paragraphs(6)
$p3.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(7)
$newPara.Range.Text = "sikapestään"
